# Updates the cryptos list (Price + Volume(1h) columns, and a couple of
# Coin/Link swaps) to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new value }
# 'D' (Price) values are forced to text with a leading apostrophe so Excel
# doesn't reinterpret strings such as '1.000' or '30.262.95' as numbers/dates,
# matching how the source data is stored (plain text) in the workbook.
$rowUpdates = [ordered]@{
    2 = @{ "D" = "30.262.95"; "E" = "  -3.31%  " }
    3 = @{ "D" = "1.928.23"; "E" = "  -3.07%  " }
    4 = @{ "D" = "1.001"; "E" = "  +0.28%  " }
    5 = @{ "D" = "246.43"; "E" = "  -2.91%  " }
    6 = @{ "D" = "0.7172"; "E" = "  -12.16%  " }
    7 = @{ "D" = "1.000"; "E" = "  +0.29%  " }
    8 = @{ "D" = "0.3267"; "E" = "  -5.69%  " }
    9 = @{ "D" = "26.48"; "E" = "  +3.03%  " }
    10 = @{ "D" = "0.06808"; "E" = "  -3.12%  " }
    11 = @{ "D" = "0.8026"; "E" = "  -4.76%  " }
    12 = @{ "D" = "0.07950"; "E" = "  -1.95%  " }
    13 = @{ "D" = "1.928.64"; "E" = "  -2.97%  " }
    14 = @{ "D" = "5.409"; "E" = "  -2.12%  " }
    15 = @{ "D" = "94.46"; "E" = "  -6.40%  " }
    16 = @{ "D" = "14.49"; "E" = "  +3.69%  " }
    17 = @{ "D" = "260.64"; "E" = "  -4.78%  " }
    18 = @{ "D" = "30.262.74"; "E" = "  -3.27%  " }
    19 = @{ "D" = "0.000007934"; "E" = "  -0.02%  " }
    20 = @{ "D" = "5.818"; "E" = "  +0.32%  " }
    21 = @{ "D" = "2.180.93"; "E" = "  -3.03%  " }
    22 = @{ "D" = "0.9999"; "E" = "  +0.16%  " }
    23 = @{ "D" = "1.000"; "E" = "  +0.27%  " }
    24 = @{ "D" = "6.876"; "E" = "  -1.31%  " }
    25 = @{ "D" = "9.670"; "E" = "  -1.15%  " }
    26 = @{ "D" = "159.89"; "E" = "  -2.65%  " }
    27 = @{ "D" = "0.1336"; "E" = "  -11.53%  " }
    28 = @{ "D" = "18.99"; "E" = "  -5.92%  " }
    29 = @{ "D" = "2.281"; "E" = "  +3.53%  " }
    30 = @{ "E" = "  +1.29%  " }
    31 = @{ "D" = "1.547"; "E" = "  -1.31%  " }
    32 = @{ "D" = "4.393"; "E" = "  -3.86%  " }
    33 = @{ "D" = "4.191"; "E" = "  -3.15%  " }
    34 = @{ "D" = "0.05075"; "E" = "  -2.09%  " }
    35 = @{ "D" = "1.195"; "E" = "  -1.54%  " }
    36 = @{ "D" = "0.7412"; "E" = "  -2.35%  " }
    37 = @{ "D" = "2.728"; "E" = "  -1.04%  " }
    38 = @{ "D" = "0.01936"; "E" = "  -3.29%  " }
    39 = @{ "E" = "  -3.29%  " }
    40 = @{ "D" = "80.03"; "E" = "  +2.30%  " }
    41 = @{ "D" = "6.571"; "E" = "  -0.64%  " }
    42 = @{ "D" = "0.4462"; "E" = "  -5.31%  " }
    43 = @{ "D" = "2.005"; "E" = "  -4.07%  " }
    44 = @{ "D" = "1.000"; "E" = "  +0.33%  " }
    45 = @{ "D" = "0.8340"; "E" = "  -2.39%  " }
    46 = @{ "D" = "102.29"; "E" = "  -2.37%  " }
    47 = @{ "D" = "9.688"; "E" = "  -2.73%  " }
    48 = @{ "D" = "7.272"; "E" = "  -3.21%  " }
    49 = @{ "D" = "36.23"; "E" = "  -1.45%  " }
    50 = @{ "B" = "Decentraland"; "C" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; "D" = "0.4108"; "E" = "  -4.96%  " }
    51 = @{ "B" = "NEARProtocol"; "C" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; "D" = "1.481"; "E" = "  +2.31%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $value = $cols[$col]
        $cellRef = "$col$row"
        if ($col -eq "D") {
            # Leading apostrophe = Excel's literal "treat as text" prefix;
            # it is not stored as part of the cell's text.
            $ws.Range($cellRef).Value = "'" + $value
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
